# Project "Sample Project" re-save: rule row 4 (row 11 on the "Rules" sheet)
# had its label cell (B11) changed from the text "R40" to the text "1".
#
# The cell keeps its existing formatting (thick-bottom-border, General
# number format) - only the stored value/type changes, from the shared
# string "R40" to a new shared string "1" (stored as literal text, not as
# a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("B11")

# Force the literal text "1" to be written as a text value (not coerced to
# the number 1) by marking the cell as Text before assigning it - mirrors
# typing into a text-formatted cell / pasting text data.
$cell.NumberFormat = "@"
$cell.Value = "1"
